$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Cells.Item(18, 1).Value = 112181755
$ws.Cells.Item(18, 2).Value = 78579
$ws.Cells.Item(18, 3).Value = '''Ovaliderad'
$ws.Cells.Item(18, 4).Value = '''NT'
$ws.Cells.Item(18, 5).Value = 2081
$ws.Cells.Item(18, 6).Value = '''Skrovellav'
$ws.Cells.Item(18, 7).Value = '''Lobaria scrobiculata'
$ws.Cells.Item(18, 8).Value = '''(Scop.) DC.'
$ws.Cells.Item(18, 9).Value = '''1'
$ws.Cells.Item(18, 16).Value = '''Vassbo, Dlr'
$ws.Cells.Item(18, 17).Value = 364894.3754826509
$ws.Cells.Item(18, 18).Value = 6872299.904270066
$ws.Cells.Item(18, 19).Value = 5
$ws.Cells.Item(18, 20).Value = '''Dalarna'
$ws.Cells.Item(18, 21).Value = '''Älvdalen'
$ws.Cells.Item(18, 22).Value = '''Dalarna'
$ws.Cells.Item(18, 23).Value = '''Idre'
$ws.Cells.Item(18, 25).Value = '''2023-06-29'
$ws.Cells.Item(18, 26).Value = '''00:00'
$ws.Cells.Item(18, 27).Value = '''2023-06-29'
$ws.Cells.Item(18, 28).Value = '''00:00'
$ws.Cells.Item(18, 29).Value = '''Påträffad under Sveaskogs naturvärdesinventering'
$ws.Cells.Item(18, 30).Value = $false
$ws.Cells.Item(18, 31).Value = $false
$ws.Cells.Item(18, 33).Value = $false
$ws.Cells.Item(18, 49).Value = '''Mimmi Persson'
$ws.Cells.Item(18, 50).Value = '''Mimmi Persson'

# Row 19
$ws.Cells.Item(19, 1).Value = 112182654
$ws.Cells.Item(19, 2).Value = 76918
$ws.Cells.Item(19, 3).Value = '''Ovaliderad'
$ws.Cells.Item(19, 4).Value = '''NT'
$ws.Cells.Item(19, 5).Value = 6437
$ws.Cells.Item(19, 6).Value = '''Blanksvart spiklav'
$ws.Cells.Item(19, 7).Value = '''Calicium denigratum'
$ws.Cells.Item(19, 8).Value = '''(Vain.) Tibell'
$ws.Cells.Item(19, 9).Value = '''1'
$ws.Cells.Item(19, 16).Value = '''Vassbo, Dlr'
$ws.Cells.Item(19, 17).Value = 364914.1726876026
$ws.Cells.Item(19, 18).Value = 6872133.350211774
$ws.Cells.Item(19, 19).Value = 5
$ws.Cells.Item(19, 20).Value = '''Dalarna'
$ws.Cells.Item(19, 21).Value = '''Älvdalen'
$ws.Cells.Item(19, 22).Value = '''Dalarna'
$ws.Cells.Item(19, 23).Value = '''Idre'
$ws.Cells.Item(19, 25).Value = '''2023-06-29'
$ws.Cells.Item(19, 26).Value = '''00:00'
$ws.Cells.Item(19, 27).Value = '''2023-06-29'
$ws.Cells.Item(19, 28).Value = '''00:00'
$ws.Cells.Item(19, 29).Value = '''Påträffad under Sveaskogs naturvärdesinventering'
$ws.Cells.Item(19, 30).Value = $false
$ws.Cells.Item(19, 31).Value = $false
$ws.Cells.Item(19, 33).Value = $false
$ws.Cells.Item(19, 49).Value = '''Mimmi Persson'
$ws.Cells.Item(19, 50).Value = '''Mimmi Persson'

# Row 20
$ws.Cells.Item(20, 1).Value = 112182890
$ws.Cells.Item(20, 2).Value = 96265
$ws.Cells.Item(20, 3).Value = '''Ovaliderad'
$ws.Cells.Item(20, 4).Value = '''LC'
$ws.Cells.Item(20, 5).Value = 219790
$ws.Cells.Item(20, 6).Value = '''Fläcknycklar'
$ws.Cells.Item(20, 7).Value = '''Dactylorhiza maculata'
$ws.Cells.Item(20, 8).Value = '''(L.) Soó'
$ws.Cells.Item(20, 9).Value = '''1'
$ws.Cells.Item(20, 16).Value = '''Vassbo, Dlr'
$ws.Cells.Item(20, 17).Value = 364947.03503229
$ws.Cells.Item(20, 18).Value = 6872307.702530573
$ws.Cells.Item(20, 19).Value = 5
$ws.Cells.Item(20, 20).Value = '''Dalarna'
$ws.Cells.Item(20, 21).Value = '''Älvdalen'
$ws.Cells.Item(20, 22).Value = '''Dalarna'
$ws.Cells.Item(20, 23).Value = '''Idre'
$ws.Cells.Item(20, 25).Value = '''2023-06-29'
$ws.Cells.Item(20, 26).Value = '''00:00'
$ws.Cells.Item(20, 27).Value = '''2023-06-29'
$ws.Cells.Item(20, 28).Value = '''00:00'
$ws.Cells.Item(20, 29).Value = '''Påträffad under Sveaskogs naturvärdesinventering'
$ws.Cells.Item(20, 30).Value = $false
$ws.Cells.Item(20, 31).Value = $false
$ws.Cells.Item(20, 33).Value = $false
$ws.Cells.Item(20, 49).Value = '''Mimmi Persson'
$ws.Cells.Item(20, 50).Value = '''Mimmi Persson'

# Row 21
$ws.Cells.Item(21, 1).Value = 112182349
$ws.Cells.Item(21, 2).Value = 77515
$ws.Cells.Item(21, 3).Value = '''Ovaliderad'
$ws.Cells.Item(21, 4).Value = '''NT'
$ws.Cells.Item(21, 5).Value = 6425
$ws.Cells.Item(21, 6).Value = '''Garnlav'
$ws.Cells.Item(21, 7).Value = '''Alectoria sarmentosa'
$ws.Cells.Item(21, 8).Value = '''(Ach.) Ach.'
$ws.Cells.Item(21, 9).Value = '''1'
$ws.Cells.Item(21, 16).Value = '''Vassbo, Dlr'
$ws.Cells.Item(21, 17).Value = 364898.4777887367
$ws.Cells.Item(21, 18).Value = 6872201.317631777
$ws.Cells.Item(21, 19).Value = 5
$ws.Cells.Item(21, 20).Value = '''Dalarna'
$ws.Cells.Item(21, 21).Value = '''Älvdalen'
$ws.Cells.Item(21, 22).Value = '''Dalarna'
$ws.Cells.Item(21, 23).Value = '''Idre'
$ws.Cells.Item(21, 25).Value = '''2023-06-29'
$ws.Cells.Item(21, 26).Value = '''00:00'
$ws.Cells.Item(21, 27).Value = '''2023-06-29'
$ws.Cells.Item(21, 28).Value = '''00:00'
$ws.Cells.Item(21, 29).Value = '''Påträffad under Sveaskogs naturvärdesinventering'
$ws.Cells.Item(21, 30).Value = $false
$ws.Cells.Item(21, 31).Value = $false
$ws.Cells.Item(21, 33).Value = $false
$ws.Cells.Item(21, 49).Value = '''Mimmi Persson'
$ws.Cells.Item(21, 50).Value = '''Mimmi Persson'

# Row 22
$ws.Cells.Item(22, 1).Value = 112182494
$ws.Cells.Item(22, 2).Value = 77267
$ws.Cells.Item(22, 3).Value = '''Ovaliderad'
$ws.Cells.Item(22, 4).Value = '''NT'
$ws.Cells.Item(22, 5).Value = 6446
$ws.Cells.Item(22, 6).Value = '''Kolflarnlav'
$ws.Cells.Item(22, 7).Value = '''Carbonicola anthracophila'
$ws.Cells.Item(22, 8).Value = '''(Nyl.) Bendiksby & Timdal'
$ws.Cells.Item(22, 9).Value = '''1'
$ws.Cells.Item(22, 16).Value = '''Vassbo, Dlr'
$ws.Cells.Item(22, 17).Value = 364938.076161085
$ws.Cells.Item(22, 18).Value = 6872236.477867194
$ws.Cells.Item(22, 19).Value = 5
$ws.Cells.Item(22, 20).Value = '''Dalarna'
$ws.Cells.Item(22, 21).Value = '''Älvdalen'
$ws.Cells.Item(22, 22).Value = '''Dalarna'
$ws.Cells.Item(22, 23).Value = '''Idre'
$ws.Cells.Item(22, 25).Value = '''2023-06-29'
$ws.Cells.Item(22, 26).Value = '''00:00'
$ws.Cells.Item(22, 27).Value = '''2023-06-29'
$ws.Cells.Item(22, 28).Value = '''00:00'
$ws.Cells.Item(22, 29).Value = '''Påträffad under Sveaskogs naturvärdesinventering'
$ws.Cells.Item(22, 30).Value = $false
$ws.Cells.Item(22, 31).Value = $false
$ws.Cells.Item(22, 33).Value = $false
$ws.Cells.Item(22, 49).Value = '''Mimmi Persson'
$ws.Cells.Item(22, 50).Value = '''Mimmi Persson'

Write-Output "Added rows 18-22"